$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 14:19"

# --- Reorder "Azerbaiyan" / "Austria" (rows 67-68) ---
# Azerbaiyan now comes first (row 67) with freshly updated stats,
# Austria moves down to row 68 keeping what used to be Azerbaiyan's slot.
$ws.Range("A67").Value = "Azerbaiyan"
$ws.Range("A68").Value = "Austria"

# --- Reorder "Montserrat" / "Islas Malvinas" (rows 214-215) ---
# Montserrat now comes first (row 214), Islas Malvinas moves to row 215.
$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# --- Updated country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7051747
$ws.Range("C4").Value = 5531
$ws.Range("D4").Value = 4300813
$ws.Range("E4").Value = 2546375
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 204559

# Row 5 - India
$ws.Range("B5").Value = 5568740
$ws.Range("C5").Value = 8635
$ws.Range("D5").Value = 4498520
$ws.Range("E5").Value = 981228
$ws.Range("G5").Value = 27
$ws.Range("H5").Value = 88992

# Row 25 - Alemania
$ws.Range("B25").Value = 275797
$ws.Range("C25").Value = 246
$ws.Range("E25").Value = 20013
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 9484

# Row 40 - Kuwait
$ws.Range("B40").Value = 100683
$ws.Range("C40").Value = 719
$ws.Range("D40").Value = 91612
$ws.Range("E40").Value = 8483
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 588

# Row 42 - Oman
$ws.Range("B42").Value = 94711
$ws.Range("C42").Value = 660
$ws.Range("D42").Value = 86195
$ws.Range("E42").Value = 7651
$ws.Range("G42").Value = 12
$ws.Range("H42").Value = 865

# Row 43 - Suecia
$ws.Range("B43").Value = 89436
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 5870

# Row 61 - Suiza
$ws.Range("B61").Value = 50664
$ws.Range("C61").Value = 286
$ws.Range("E61").Value = 6811
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 2053

# Row 67 - now Azerbaiyan (fresh figures)
$ws.Range("B67").Value = 39378
$ws.Range("C67").Value = 98
$ws.Range("D67").Value = 36949
$ws.Range("E67").Value = 1851
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 578

# Row 68 - now Austria (carries forward the previous Austria figures)
$ws.Range("B68").Value = 39303
$ws.Range("C68").Value = 645
$ws.Range("D68").Value = 30312
$ws.Range("E68").Value = 8220
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 771

# Row 71 - Estado de Palestina
$ws.Range("B71").Value = 36580
$ws.Range("C71").Value = 429
$ws.Range("D71").Value = 25469
$ws.Range("E71").Value = 10842
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 269

# Row 76 - Libia
$ws.Range("B76").Value = 29446
$ws.Range("C76").Value = 650
$ws.Range("D76").Value = 15913
$ws.Range("E76").Value = 13073
$ws.Range("G76").Value = 10
$ws.Range("H76").Value = 460

# Row 80 - Dinamarca
$ws.Range("B80").Value = 23799
$ws.Range("C80").Value = 476
$ws.Range("D80").Value = 18035
$ws.Range("E80").Value = 5123
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 641

# Row 106 - Haiti
$ws.Range("B106").Value = 8633
$ws.Range("C106").Value = 9
$ws.Range("E106").Value = 1928
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 223

# Row 214 - now Montserrat (D/H swap with old row 215)
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215 - now Islas Malvinas (D/H swap with old row 214)
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
